$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Functionality" table (Table3) with 5 new TODO rows and
# re-apply it so the AutoFilter range/table range cover B2:C20.
$lo = $ws.ListObjects("Table3")
$lo.Resize($ws.Range("B2:C20")) | Out-Null

$ws.Range("B16").Value = "Cache base stats on load"
$ws.Range("C16").Value = "todo"

$ws.Range("B17").Value = "cache equipped stats on load"
$ws.Range("C17").Value = "todo"

$ws.Range("B18").Value = "cache traits on load"
$ws.Range("C18").Value = "todo"

$ws.Range("B19").Value = "update equipped stats on gear change"
$ws.Range("C19").Value = "todo"

$ws.Range("B20").Value = "update cached traits on trait update"
$ws.Range("C20").Value = "todo"

# Filter the "Column1" (status) column down to rows marked "todo" only,
# hiding the already-"done" rows.
$lo.Range.AutoFilter(2, @("todo")) | Out-Null

# Move the active selection to where the newly-added rows now start.
$ws.Range("B21").Select() | Out-Null
